# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7571
$ws1.Range("F3").Value = 92
$ws1.Range("F4").Value = 72
$ws1.Range("F5").Value = 4587
$ws1.Range("F7").Value = 588
$ws1.Range("F8").Value = 608
$ws1.Range("F9").Value = 441
$ws1.Range("F11").Value = 419
$ws1.Range("F12").Value = 756
$ws1.Range("F13").Value = 24
$ws1.Range("F14").Value = 62
$ws1.Range("F15").Value = 246
$ws1.Range("F16").Value = 13
$ws1.Range("F17").Value = 245
$ws1.Range("F18").Value = 129
$ws1.Range("F19").Value = 380
$ws1.Range("F23").Value = 547
$ws1.Range("F24").Value = 2153
$ws1.Range("F25").Value = 678
$ws1.Range("F26").Value = 36
$ws1.Range("F27").Value = 34
$ws1.Range("F30").Value = 38

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 286

# --- 本地生活 (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 427

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 427
$ws4.Range("F3").Value = 7571
$ws4.Range("F4").Value = 92
$ws4.Range("F5").Value = 72
$ws4.Range("F6").Value = 286
$ws4.Range("F7").Value = 4588
$ws4.Range("F9").Value = 588
$ws4.Range("F10").Value = 608
$ws4.Range("F11").Value = 441
$ws4.Range("F14").Value = 419
$ws4.Range("F18").Value = 756
$ws4.Range("F19").Value = 24
$ws4.Range("F21").Value = 246
$ws4.Range("F23").Value = 13
$ws4.Range("F26").Value = 245
$ws4.Range("F27").Value = 129
$ws4.Range("F28").Value = 380
$ws4.Range("F32").Value = 547
$ws4.Range("F33").Value = 2153
$ws4.Range("F34").Value = 678
$ws4.Range("F35").Value = 36
$ws4.Range("F36").Value = 34
$ws4.Range("F39").Value = 38
